$d = $word.ActiveDocument

# Add a black, 1/4pt (sz=4), single-line border to every edge of every
# table in the document: outside box (top/left/bottom/right) plus the
# inside horizontal/vertical gridlines.
#
# wdLineStyleSingle   = 1
# wdLineWidth025pt    = 2   (-> w:sz="4", eighths of a point)
# wdColorAutomatic/Black (index 0 -> w:color="000000")
#
# wdBorderTop = -1, wdBorderLeft = -2, wdBorderBottom = -3,
# wdBorderRight = -4, wdBorderHorizontal = -5, wdBorderVertical = -6

foreach ($table in $d.Tables) {
    for ($i = -1; $i -ge -6; $i--) {
        $border = $table.Borders.Item($i)
        $border.LineStyle = 1
        $border.LineWidth = 2
        $border.Color = 0
    }
}
